$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Apply formatting to new cells by copying format from existing styled neighbors ---
$ws.Range("S1").Copy()
$ws.Range("T1").PasteSpecial(-4122)
$ws.Range("A7").Copy()
$ws.Range("A8:A11").PasteSpecial(-4122)

# --- New-string priming in target insertion order: OffsetF, OffsetA, RD Single, TD Single, 1Pair-B ---
$ws.Cells.Item(7, 2).Value = "OffsetF"
$ws.Cells.Item(8, 2).Value = "OffsetA"
$ws.Cells.Item(9, 2).Value = "RD Single"
$ws.Cells.Item(10, 2).Value = "TD Single"
$ws.Cells.Item(2, 12).Value = "1Pair-B"

# --- Row 1 header (unchanged values, only T1 is new) ---
$ws.Range("T1").Value = 18

# --- Row 2 ---
$ws.Cells.Item(2, 1).Value = 0
$ws.Cells.Item(2, 2).Value = "HKL"
$ws.Cells.Item(2, 3).Value = "[1, 1, 0]"
$ws.Cells.Item(2, 4).Value = "[2, 0, 0]"
$ws.Cells.Item(2, 5).Value = "[2, 1, 1]"
$ws.Cells.Item(2, 6).Value = "[2, 2, 0]"
$ws.Cells.Item(2, 7).Value = "[3, 1, 0]"
$ws.Cells.Item(2, 8).Value = "[2, 2, 2]"
$ws.Cells.Item(2, 9).Value = "[3, 2, 1]"
$ws.Cells.Item(2, 10).Value = "[4, 0, 0]"
$ws.Cells.Item(2, 11).Value = "1Pair-A"
$ws.Cells.Item(2, 12).Value = "1Pair-B"
$ws.Cells.Item(2, 13).Value = "2Pairs-A"
$ws.Cells.Item(2, 14).Value = "2Pairs-B"
$ws.Cells.Item(2, 15).Value = "3Pairs-A"
$ws.Cells.Item(2, 16).Value = "3Pairs-B"
$ws.Cells.Item(2, 17).Value = "3Pairs-C"
$ws.Cells.Item(2, 18).Value = "4Pairs"
$ws.Cells.Item(2, 19).Value = "5A4F"
$ws.Cells.Item(2, 20).Value = "MaxUnique"

# --- Row 3 ---
$ws.Cells.Item(3, 1).Value = 1
$ws.Cells.Item(3, 2).Value = "Equal Angle"
$ws.Cells.Item(3, 3).Value = 0.9766066282420749
$ws.Cells.Item(3, 4).Value = 0.806621037463977
$ws.Cells.Item(3, 5).Value = 1.07628242074928
$ws.Cells.Item(3, 6).Value = 0.9766066282420749
$ws.Cells.Item(3, 7).Value = 0.8578314121037464
$ws.Cells.Item(3, 8).Value = 1.262017291066282
$ws.Cells.Item(3, 9).Value = 1.05271613832853
$ws.Cells.Item(3, 10).Value = 0.806621037463977
$ws.Cells.Item(3, 11).Value = 0.9766066282420749
$ws.Cells.Item(3, 12).Value = 1.07628242074928
$ws.Cells.Item(3, 13).Value = 0.9414517291066282
$ws.Cells.Item(3, 14).Value = 0.9414517291066282
$ws.Cells.Item(3, 15).Value = 0.9135782901056676
$ws.Cells.Item(3, 16).Value = 0.9531700288184437
$ws.Cells.Item(3, 17).Value = 0.9531700288184437
$ws.Cells.Item(3, 18).Value = 0.9590291786743514
$ws.Cells.Item(3, 19).Value = 0.9590291786743514
$ws.Cells.Item(3, 20).Value = 1.005345821325648

# --- Row 4 ---
$ws.Cells.Item(4, 1).Value = 2
$ws.Cells.Item(4, 2).Value = "CLR"
$ws.Cells.Item(4, 3).Value = 1.009662181161954
$ws.Cells.Item(4, 4).Value = 0.9665531897185299
$ws.Cells.Item(4, 5).Value = 0.9979974765673288
$ws.Cells.Item(4, 6).Value = 1.009662181161954
$ws.Cells.Item(4, 7).Value = 0.9773453179064604
$ws.Cells.Item(4, 8).Value = 1.002136977492128
$ws.Cells.Item(4, 9).Value = 1.001787173412981
$ws.Cells.Item(4, 10).Value = 0.9665531897185299
$ws.Cells.Item(4, 11).Value = 1.009662181161954
$ws.Cells.Item(4, 12).Value = 0.9979974765673288
$ws.Cells.Item(4, 13).Value = 0.9822753331429294
$ws.Cells.Item(4, 14).Value = 0.9822753331429294
$ws.Cells.Item(4, 15).Value = 0.9806319947307731
$ws.Cells.Item(4, 16).Value = 0.9914042824826043
$ws.Cells.Item(4, 17).Value = 0.9914042824826043
$ws.Cells.Item(4, 18).Value = 0.9959687571524418
$ws.Cells.Item(4, 19).Value = 0.9959687571524418
$ws.Cells.Item(4, 20).Value = 0.9925803860432304

# --- Row 5 ---
$ws.Cells.Item(5, 1).Value = 3
$ws.Cells.Item(5, 2).Value = "BT8Hex"
$ws.Cells.Item(5, 3).Value = 1.025020716091789
$ws.Cells.Item(5, 4).Value = 0.9431103096578571
$ws.Cells.Item(5, 5).Value = 0.9998160823045482
$ws.Cells.Item(5, 6).Value = 1.025020716091789
$ws.Cells.Item(5, 7).Value = 0.9655641955633906
$ws.Cells.Item(5, 8).Value = 1.006412446476209
$ws.Cells.Item(5, 9).Value = 1.007384609498003
$ws.Cells.Item(5, 10).Value = 0.9431103096578571
$ws.Cells.Item(5, 11).Value = 1.025020716091789
$ws.Cells.Item(5, 12).Value = 0.9998160823045482
$ws.Cells.Item(5, 13).Value = 0.9714631959812027
$ws.Cells.Item(5, 14).Value = 0.9714631959812027
$ws.Cells.Item(5, 15).Value = 0.9694968625085987
$ws.Cells.Item(5, 16).Value = 0.9893157026847312
$ws.Cells.Item(5, 17).Value = 0.9893157026847312
$ws.Cells.Item(5, 18).Value = 0.9982419560364956
$ws.Cells.Item(5, 19).Value = 0.9982419560364956
$ws.Cells.Item(5, 20).Value = 0.991218059931966

# --- Row 6 ---
$ws.Cells.Item(6, 1).Value = 4
$ws.Cells.Item(6, 2).Value = "Spiral"
$ws.Cells.Item(6, 3).Value = 0.9906822651826908
$ws.Cells.Item(6, 4).Value = 0.9934232063263552
$ws.Cells.Item(6, 5).Value = 0.9971114466606487
$ws.Cells.Item(6, 6).Value = 0.9906822651826908
$ws.Cells.Item(6, 7).Value = 0.9895999663079804
$ws.Cells.Item(6, 8).Value = 1.002970065815645
$ws.Cells.Item(6, 9).Value = 0.9959255459470567
$ws.Cells.Item(6, 10).Value = 0.9934232063263552
$ws.Cells.Item(6, 11).Value = 0.9906822651826908
$ws.Cells.Item(6, 12).Value = 0.9971114466606487
$ws.Cells.Item(6, 13).Value = 0.9952673264935019
$ws.Cells.Item(6, 14).Value = 0.9952673264935019
$ws.Cells.Item(6, 15).Value = 0.9933782064316614
$ws.Cells.Item(6, 16).Value = 0.9937389727232316
$ws.Cells.Item(6, 17).Value = 0.9937389727232316
$ws.Cells.Item(6, 18).Value = 0.9929747958380963
$ws.Cells.Item(6, 19).Value = 0.9929747958380963
$ws.Cells.Item(6, 20).Value = 0.9949520827067295

# --- Row 7 ---
$ws.Cells.Item(7, 1).Value = 5
$ws.Cells.Item(7, 2).Value = "OffsetF"
$ws.Cells.Item(7, 3).Value = 1.120092255082555
$ws.Cells.Item(7, 4).Value = 1.046331908572508
$ws.Cells.Item(7, 5).Value = 0.9000153727474539
$ws.Cells.Item(7, 6).Value = 1.120092255082555
$ws.Cells.Item(7, 7).Value = 1.095417664142477
$ws.Cells.Item(7, 8).Value = 0.6638777778453894
$ws.Cells.Item(7, 9).Value = 0.9557535756961746
$ws.Cells.Item(7, 10).Value = 1.046331908572508
$ws.Cells.Item(7, 11).Value = 1.120092255082555
$ws.Cells.Item(7, 12).Value = 0.9000153727474539
$ws.Cells.Item(7, 13).Value = 0.9731736406599811
$ws.Cells.Item(7, 14).Value = 0.9731736406599811
$ws.Cells.Item(7, 15).Value = 1.01392164848748
$ws.Cells.Item(7, 16).Value = 1.022146512134173
$ws.Cells.Item(7, 17).Value = 1.022146512134173
$ws.Cells.Item(7, 18).Value = 1.046632947871268
$ws.Cells.Item(7, 19).Value = 1.046632947871268
$ws.Cells.Item(7, 20).Value = 0.9635814256810932

# --- Row 8 ---
$ws.Cells.Item(8, 1).Value = 6
$ws.Cells.Item(8, 2).Value = "OffsetA"
$ws.Cells.Item(8, 3).Value = 0.9196222567790814
$ws.Cells.Item(8, 4).Value = 1.123183439351588
$ws.Cells.Item(8, 5).Value = 1.003424219258546
$ws.Cells.Item(8, 6).Value = 0.9196222567790814
$ws.Cells.Item(8, 7).Value = 1.018026418624242
$ws.Cells.Item(8, 8).Value = 1.047367271581791
$ws.Cells.Item(8, 9).Value = 0.9761861983916611
$ws.Cells.Item(8, 10).Value = 1.123183439351588
$ws.Cells.Item(8, 11).Value = 0.9196222567790814
$ws.Cells.Item(8, 12).Value = 1.003424219258546
$ws.Cells.Item(8, 13).Value = 1.063303829305067
$ws.Cells.Item(8, 14).Value = 1.063303829305067
$ws.Cells.Item(8, 15).Value = 1.048211359078125
$ws.Cells.Item(8, 16).Value = 1.015409971796405
$ws.Cells.Item(8, 17).Value = 1.015409971796405
$ws.Cells.Item(8, 18).Value = 0.9914630430420742
$ws.Cells.Item(8, 19).Value = 0.9914630430420742
$ws.Cells.Item(8, 20).Value = 1.014634967331151

# --- Row 9 ---
$ws.Cells.Item(9, 1).Value = 7
$ws.Cells.Item(9, 2).Value = "RD Single"
$ws.Cells.Item(9, 3).Value = 1.98
$ws.Cells.Item(9, 4).Value = 0.21
$ws.Cells.Item(9, 5).Value = 0.84
$ws.Cells.Item(9, 6).Value = 1.98
$ws.Cells.Item(9, 7).Value = 0.64
$ws.Cells.Item(9, 8).Value = 0.6899999999999999
$ws.Cells.Item(9, 9).Value = 1.14
$ws.Cells.Item(9, 10).Value = 0.21
$ws.Cells.Item(9, 11).Value = 1.98
$ws.Cells.Item(9, 12).Value = 0.84
$ws.Cells.Item(9, 13).Value = 0.525
$ws.Cells.Item(9, 14).Value = 0.525
$ws.Cells.Item(9, 15).Value = 0.5633333333333334
$ws.Cells.Item(9, 16).Value = 1.01
$ws.Cells.Item(9, 17).Value = 1.01
$ws.Cells.Item(9, 18).Value = 1.2525
$ws.Cells.Item(9, 19).Value = 1.2525
$ws.Cells.Item(9, 20).Value = 0.9166666666666665

# --- Row 10 ---
$ws.Cells.Item(10, 1).Value = 8
$ws.Cells.Item(10, 2).Value = "TD Single"
$ws.Cells.Item(10, 3).Value = 0.84
$ws.Cells.Item(10, 4).Value = 0.44
$ws.Cells.Item(10, 5).Value = 1.27
$ws.Cells.Item(10, 6).Value = 0.84
$ws.Cells.Item(10, 7).Value = 0.68
$ws.Cells.Item(10, 8).Value = 1.52
$ws.Cells.Item(10, 9).Value = 1.15
$ws.Cells.Item(10, 10).Value = 0.44
$ws.Cells.Item(10, 11).Value = 0.84
$ws.Cells.Item(10, 12).Value = 1.27
$ws.Cells.Item(10, 13).Value = 0.855
$ws.Cells.Item(10, 14).Value = 0.855
$ws.Cells.Item(10, 15).Value = 0.7966666666666667
$ws.Cells.Item(10, 16).Value = 0.85
$ws.Cells.Item(10, 17).Value = 0.85
$ws.Cells.Item(10, 18).Value = 0.8474999999999999
$ws.Cells.Item(10, 19).Value = 0.8474999999999999
$ws.Cells.Item(10, 20).Value = 0.9833333333333334

# --- Row 11 ---
$ws.Cells.Item(11, 1).Value = 9
$ws.Cells.Item(11, 2).Value = "HexGrid-90degTilt5degRes"
$ws.Cells.Item(11, 3).Value = 0.9950810218029537
$ws.Cells.Item(11, 4).Value = 0.9948803990747337
$ws.Cells.Item(11, 5).Value = 0.9943622104465786
$ws.Cells.Item(11, 6).Value = 0.9950810218029537
$ws.Cells.Item(11, 7).Value = 0.9916272098834858
$ws.Cells.Item(11, 8).Value = 0.9951506684726222
$ws.Cells.Item(11, 9).Value = 0.99493714788777
$ws.Cells.Item(11, 10).Value = 0.9948803990747337
$ws.Cells.Item(11, 11).Value = 0.9950810218029537
$ws.Cells.Item(11, 12).Value = 0.9943622104465786
$ws.Cells.Item(11, 13).Value = 0.9946213047606561
$ws.Cells.Item(11, 14).Value = 0.9946213047606561
$ws.Cells.Item(11, 15).Value = 0.9936232731349327
$ws.Cells.Item(11, 16).Value = 0.9947745437747554
$ws.Cells.Item(11, 17).Value = 0.9947745437747554
$ws.Cells.Item(11, 18).Value = 0.9948511632818049
$ws.Cells.Item(11, 19).Value = 0.9948511632818049
$ws.Cells.Item(11, 20).Value = 0.9943397762613574
